$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Update column C (rows 2-27) date serial value from 45317 to 45318
$ws.Range("C2:C27").Value = 45318
